$d = $word.ActiveDocument
$wdYellow = 7

function Find-ParagraphByText($doc, [string]$needle) {
    $n = $doc.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $needle) {
            return $p
        }
    }
    throw "Paragraph not found: $needle"
}

# ---------------------------------------------------------------------------
# Simple highlight-only paragraphs (single or multi-run, highlight every run)
# ---------------------------------------------------------------------------
$simpleHighlightTexts = @(
    "Active Alpha Subunit shall seek Adenylyl Cyclase.",
    "When GTP binds with the activated GPC Receptor, its GDP shall break off.",
    "When the Trimeric G-Protein separates from the GPC Receptor, the alpha subunit shall split off from the Beta-Gamma Complex.",
    "The Alpha Subunit, once separated from the Beta-Gamma Complex shall seek an Adenylyl Cyclase with which to bind.",
    "The Alpha Subunit shall bind with the Adenylyl Cyclase.",
    "BioRube Bot shall allow the user to spawn an Adenylyl Cyclase on the Cell Membrane.",
    "The Adenylyl Cyclase shall await activation via the Alpha-Beta Subunit.",
    "Once activated the Adenylyl Cyclase shall undergo a transformation.",
    "Checkbox two shall become checked after the Trimeric G-Protein binds with the G-Protein Coupled Receptor.",
    "Checkbox three shall become checked when the Trimeric G-Protein binds with a GTP and breaks apart.",
    "Checkbox four shall become checked when the alpha subunit binds with the Adenylyl Cyclase."
)

foreach ($t in $simpleHighlightTexts) {
    $p = Find-ParagraphByText $d $t
    $p.Range.HighlightColorIndex = $wdYellow
}

# ---------------------------------------------------------------------------
# Paragraph: "The Adenylyl Cyclase shall adhere to the cell as depicted in
# Figure 1." -- add highlight to every run, and consolidate the REF field's
# instrText runs into a single run (as Word does on a field-code update).
# ---------------------------------------------------------------------------
$pAdhere = Find-ParagraphByText $d "The Adenylyl Cyclase shall adhere to the cell as depicted in Figure 1."
$adhereBody = ""
$adhereBody += '<w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve">The Adenylyl Cyclase shall adhere to the cell as depicted in </w:t></w:r>'
$adhereBody += '<w:r><w:rPr><w:b/><w:bCs/><w:highlight w:val="yellow"/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r>'
$adhereBody += '<w:r><w:rPr><w:b/><w:bCs/><w:highlight w:val="yellow"/></w:rPr><w:instrText xml:space="preserve"> REF _Ref86479136 \h  \* MERGEFORMAT </w:instrText></w:r>'
$adhereBody += '<w:r><w:rPr><w:b/><w:bCs/><w:highlight w:val="yellow"/></w:rPr></w:r>'
$adhereBody += '<w:r><w:rPr><w:b/><w:bCs/><w:highlight w:val="yellow"/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r>'
$adhereBody += '<w:r><w:rPr><w:b/><w:bCs/><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve">Figure </w:t></w:r>'
$adhereBody += '<w:r><w:rPr><w:b/><w:bCs/><w:noProof/><w:highlight w:val="yellow"/></w:rPr><w:t>1</w:t></w:r>'
$adhereBody += '<w:r><w:rPr><w:b/><w:bCs/><w:highlight w:val="yellow"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r>'
$adhereBody += '<w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>.</w:t></w:r>'
$adhereXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="434E389F" w14:textId="0CBB1E36" w:rsidR="00504533" w:rsidRDefault="00504533" w:rsidP="00504533">' + $adhereBody + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pAdhere.Range.InsertXML($adhereXml)

# ---------------------------------------------------------------------------
# Paragraph: "Figure 1: Adenylyl Cyclase Cell Binding" -- turn the fldSimple
# SEQ field into an explicit begin/instrText/separate/end field (no highlight
# here, this paragraph is not highlighted in the diff).
# ---------------------------------------------------------------------------
$pFigure = Find-ParagraphByText $d "Figure 1: Adenylyl Cyclase Cell Binding"
$figureBody = ""
$figureBody += '<w:pPr><w:pStyle w:val="Caption"/><w:jc w:val="center"/></w:pPr>'
$figureBody += '<w:bookmarkStart w:id="0" w:name="_Ref86479136"/>'
$figureBody += '<w:r><w:t xml:space="preserve">Figure </w:t></w:r>'
$figureBody += '<w:r><w:fldChar w:fldCharType="begin"/></w:r>'
$figureBody += '<w:r><w:instrText xml:space="preserve"> SEQ Figure \* ARABIC </w:instrText></w:r>'
$figureBody += '<w:r><w:fldChar w:fldCharType="separate"/></w:r>'
$figureBody += '<w:r><w:rPr><w:noProof/></w:rPr><w:t>1</w:t></w:r>'
$figureBody += '<w:r><w:rPr><w:noProof/></w:rPr><w:fldChar w:fldCharType="end"/></w:r>'
$figureBody += '<w:bookmarkEnd w:id="0"/>'
$figureBody += '<w:r><w:t>: Adenylyl Cyclase Cell Binding</w:t></w:r>'
$figureXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="29D6BBD7" w14:textId="0182C088" w:rsidR="00504533" w:rsidRDefault="00504533" w:rsidP="00504533">' + $figureBody + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pFigure.Range.InsertXML($figureXml)

# ---------------------------------------------------------------------------
# Paragraph: "Once the Alpha Subunit binds with the Adenylyl Cyclase, the
# Adenylyl Cyclase shall become active." -- split into two runs, highlight
# each run, and highlight the paragraph mark (pPr/rPr) too.
# ---------------------------------------------------------------------------
$pActive = Find-ParagraphByText $d "Once the Alpha Subunit binds with the Adenylyl Cyclase, the Adenylyl Cyclase shall become active."
$activeBody = ""
$activeBody += '<w:pPr><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr>'
$activeBody += '<w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve">Once the Alpha Subunit binds with the Adenylyl Cyclase, the Adenylyl Cyclase shall become </w:t></w:r>'
$activeBody += '<w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>active.</w:t></w:r>'
$activeXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="2012763F" w14:textId="51C396E2" w:rsidR="00504533" w:rsidRDefault="00504533" w:rsidP="00504533">' + $activeBody + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pActive.Range.InsertXML($activeXml)

Write-Host "Edits applied."
